$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell $ws "D2" '59.319.28'
Set-TextCell $ws "E2" '  +1.17%  '
Set-TextCell $ws "D3" '2.549.65'
Set-TextCell $ws "E3" '  +4.97%  '
Set-TextCell $ws "D4" '1.00'
Set-TextCell $ws "E4" '  +0.17%  '
Set-TextCell $ws "D5" '541.82'
Set-TextCell $ws "E5" '  +1.64%  '
Set-TextCell $ws "D6" '144.76'
Set-TextCell $ws "E6" '  +0.18%  '
Set-TextCell $ws "D7" '1.00'
Set-TextCell $ws "E7" '  +0.06%  '
Set-TextCell $ws "E8" '  +0.68%  '
Set-TextCell $ws "D9" '2.548.06'
Set-TextCell $ws "E9" '  +4.53%  '
Set-TextCell $ws "E10" '  +2.50%  '
Set-TextCell $ws "E11" '  +0.85%  '
Set-TextCell $ws "D12" '5.60'
Set-TextCell $ws "E12" '  +5.98%  '
Set-TextCell $ws "D13" '0.357'
Set-TextCell $ws "E13" '  +3.09%  '
Set-TextCell $ws "D14" '2.996.82'
Set-TextCell $ws "E14" '  +4.44%  '
Set-TextCell $ws "D15" '23.83'
Set-TextCell $ws "E15" '  +0.64%  '
Set-TextCell $ws "D16" '59.538.23'
Set-TextCell $ws "E16" '  +1.71%  '
Set-TextCell $ws "E17" '  +3.54%  '
Set-TextCell $ws "D18" '2.545.62'
Set-TextCell $ws "E18" '  +2.27%  '
Set-TextCell $ws "D19" '11.24'
Set-TextCell $ws "E19" '  +2.30%  '
Set-TextCell $ws "E20" '  +0.05%  '
Set-TextCell $ws "D21" '326.45'
Set-TextCell $ws "E21" '  +1.48%  '
Set-TextCell $ws "D22" '1.00'
Set-TextCell $ws "E22" '  +3.67%  '
Set-TextCell $ws "E23" '  +3.81%  '
Set-TextCell $ws "D24" '62.38'
Set-TextCell $ws "E24" '  +3.48%  '
Set-TextCell $ws "D25" '0.438'
Set-TextCell $ws "E25" '  -1.62%  '
Set-TextCell $ws "E26" '  +4.03%  '
Set-TextCell $ws "D27" '0.993'
Set-TextCell $ws "E27" '  +1.93%  '
Set-TextCell $ws "D28" '8.07'
Set-TextCell $ws "E28" '  +6.29%  '
Set-TextCell $ws "D29" '6.82'
Set-TextCell $ws "E29" '  +3.91%  '
Set-TextCell $ws "D30" '0.0₃0784'
Set-TextCell $ws "E30" '  +3.43%  '
Set-TextCell $ws "E31" '  +1.84%  '
Set-TextCell $ws "D32" '1.23'
Set-TextCell $ws "E32" '  +2.94%  '
Set-TextCell $ws "E33" '  +11.26%  '
Set-TextCell $ws "D34" '0.999'
Set-TextCell $ws "E34" '  +0.05%  '
Set-TextCell $ws "D35" '157.60'
Set-TextCell $ws "E35" '  +1.53%  '
Set-TextCell $ws "E36" '  +1.77%  '
Set-TextCell $ws "E37" '  +0.81%  '
Set-TextCell $ws "E38" '  -2.23%  '
Set-TextCell $ws "D39" '5.63'
Set-TextCell $ws "E39" '  -1.01%  '
Set-TextCell $ws "D40" '36.96'
Set-TextCell $ws "E40" '  +2.30%  '
Set-TextCell $ws "D41" '299.55'
Set-TextCell $ws "E41" '  -2.64%  '
Set-TextCell $ws "D42" '3.72'
Set-TextCell $ws "E42" '  +1.47%  '
Set-TextCell $ws "D43" '0.831'
Set-TextCell $ws "E43" '  +1.20%  '
Set-TextCell $ws "D44" '0.997'
Set-TextCell $ws "E44" '  -0.10%  '
Set-TextCell $ws "D45" '0.604'
Set-TextCell $ws "E45" '  +4.83%  '
Set-TextCell $ws "D46" '10.80'
Set-TextCell $ws "E46" '  +0.69%  '
Set-TextCell $ws "D47" '0.0937'
Set-TextCell $ws "E47" '  +0.48%  '
Set-TextCell $ws "D48" '18.86'
Set-TextCell $ws "E48" '  +4.05%  '
Set-TextCell $ws "D49" '122.58'
Set-TextCell $ws "E49" '  +1.19%  '
Set-TextCell $ws "E50" '  -0.63%  '
Set-TextCell $ws "D51" '0.0228'
Set-TextCell $ws "E51" '  +0.80%  '
